$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after the current last row (71), shifting
# cells down. Inserting via Range.Insert copies the formatting of the
# row immediately above into the freshly inserted rows, which matches
# the s="9" (aula/sessão) and s="1" (nome da aula/observação) styles
# used by row 71.
$ws.Range("A72:D73").Insert(-4121)

# Row 72: "69. Executando a aplicação via .jar" / empacotar em .jar via maven
$ws.Range("A72").Value = 69
$ws.Range("B72").Value = "14. Final"
$ws.Range("C72").Value = "`r`n69. Executando a aplicação via .jar"

# Row 73: same aula/sessão/nome da aula, different observação about the "/" redirect caveat
$ws.Range("A73").Value = 69
$ws.Range("B73").Value = "14. Final"
$ws.Range("C73").Value = "`r`n69. Executando a aplicação via .jar"

# Fill the "observação" column last, row 73 before row 72, so that the new
# shared-string entries land in the same order as the target workbook
# (index 129 = row73's text, index 130 = row72's text).
$ws.Range("D73").Value = "6:43`r`nem ações que não são do tipo `"redirect`" é preciso retirar a barra `"/`" no inicio de cada link existente na aplicação pois gera erro/conflito de path entre o Thymeleaf e o Spring boot ao executar o app diretamente pelo .jar"
$ws.Range("D72").Value = "para executar a aplicação fora da IDE é necessário empacotar as classes executando e fazendo build através do maven e em `"goals`" adicionar o `"package -e`" para que seja criado o arquivo .jar na pasta target da aplicação. Após o arquivo criado, basta executalo pelo cmd atraves do comando `"java -jar nomeDaSuaAplicacao.jar`" (sem aspas)"

# Row heights from the diff
$ws.Rows.Item(72).RowHeight = 30
$ws.Rows.Item(73).RowHeight = 60

# Grow the Excel table (Tabela1) to include the two new rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D73"))

# Match the new active selection reported in the diff
$ws.Range("D73").Select() | Out-Null
